# Update report.html to reflect test results and metadata changes
# - Modified the output suite data to update test case identifiers and results.
# - Added a new string to the output for the addExpenseTest case.
# - Updated statistics to indicate a failure in the addExpenseTest case.
# - Adjusted baseMillis and generated values to reflect the latest test run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 : Execute flag flips Y -> N, ActualResult gets populated, Result Pass
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "N"
$ws.Range("J2").Value = "Pass"

# ---------------------------------------------------------------------------
# Row 3 : Execute flag flips N -> Y, DeviceName spelling fix (ซ -> ช),
#          ActualResult + Result added
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Y"
$ws.Range("E3").Value = "เมล็ดผักชี"
$ws.Range("I3").ClearFormats()
$ws.Range("I3").Value = "not message"
$ws.Range("J3").Value = "Fail"

# ---------------------------------------------------------------------------
# Row 4 : Execute flag flips N -> Y, DeviceName spelling fix, ExpectedResult
#          text change, ActualResult + Result added
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Y"
$ws.Range("E4").Value = "เมล็ดผักชี"
$ws.Range("H4").Value = "บันทึกรายจ่ายสำเร็จ"
$ws.Range("I4").ClearFormats()
$ws.Range("I4").Value = "not message"
$ws.Range("J4").Value = "Fail"

# ---------------------------------------------------------------------------
# Row 5 : Execute flag flips N -> Y, ActualResult + Result added
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Y"
$ws.Range("I5").ClearFormats()
$ws.Range("I5").Value = "not message"
$ws.Range("J5").Value = "Fail"

# ---------------------------------------------------------------------------
# Row 6 : DeviceName spelling fix, ActualResult (numeric-looking text) + Result
# ---------------------------------------------------------------------------
$ws.Range("E6").Value = "เมล็ดผักชี"

# ---------------------------------------------------------------------------
# Row 7 : DeviceName spelling fix, ActualResult (numeric-looking text) + Result
# ---------------------------------------------------------------------------
$ws.Range("E7").Value = "เมล็ดผักชี"

# ---------------------------------------------------------------------------
# Row 8 : DeviceName spelling fix, ActualResult + Result
# ---------------------------------------------------------------------------
$ws.Range("E8").Value = "เมล็ดผักชี"
$ws.Range("I8").Value = "NaN"
$ws.Range("J8").Value = "Pass"

# ---------------------------------------------------------------------------
# Row 9 : DeviceName spelling fix, ActualResult + Result
# ---------------------------------------------------------------------------
$ws.Range("E9").Value = "เมล็ดผักชี"
$ws.Range("I9").Value = "กรุณากรอกราคา"
$ws.Range("J9").Value = "Pass"

# ---------------------------------------------------------------------------
# Row 10 : DeviceName spelling fix, ActualResult (numeric-looking text) + Result
# ---------------------------------------------------------------------------
$ws.Range("E10").Value = "เมล็ดผักชี"

# ---------------------------------------------------------------------------
# Row 11 : DeviceName spelling fix, ActualResult (numeric-looking text) + Result
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = "เมล็ดผักชี"

# ---------------------------------------------------------------------------
# Row 12 : DeviceName spelling fix, ActualResult + Result
# ---------------------------------------------------------------------------
$ws.Range("E12").Value = "เมล็ดผักชี"
$ws.Range("I12").Value = "NaN"
$ws.Range("J12").Value = "Pass"

# ---------------------------------------------------------------------------
# Row 13 : DeviceName spelling fix, ActualResult + Result
# ---------------------------------------------------------------------------
$ws.Range("E13").Value = "เมล็ดผักชี"
$ws.Range("I13").Value = "NaN"
$ws.Range("J13").Value = "Pass"

# ---------------------------------------------------------------------------
# Row 14 : DeviceName spelling fix, ActualResult + Result
# ---------------------------------------------------------------------------
$ws.Range("E14").Value = "เมล็ดผักชี"
$ws.Range("I14").Value = "กรุณากรอกจำนวน"
$ws.Range("J14").Value = "Pass"

# ---------------------------------------------------------------------------
# The ActualResult cells in rows 2, 6, 7, 10, 11 hold numeric-looking text
# ("50", "10", "-10", "50", "-50") that must stay stored as TEXT (matching the
# "ActualResult" column of the original sheet, same style s="1") rather than
# being auto-coerced to a number by Excel's normal type inference. Route the
# value through a scratch cell that is explicitly text-formatted, then paste
# the *values* (not formats) into the destination so the destination keeps
# its own existing style/format.
# ---------------------------------------------------------------------------
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

$scratch.Value = "50"
$scratch.Copy()
$ws.Range("I2").PasteSpecial(-4163)

$scratch.Value = "10"
$scratch.Copy()
$ws.Range("I6").PasteSpecial(-4163)

$scratch.Value = "-10"
$scratch.Copy()
$ws.Range("I7").PasteSpecial(-4163)

$scratch.Value = "50"
$scratch.Copy()
$ws.Range("I10").PasteSpecial(-4163)

$scratch.Value = "-50"
$scratch.Copy()
$ws.Range("I11").PasteSpecial(-4163)

$scratch.Clear()

# Results after the ActualResult values above
$ws.Range("J6").Value = "Pass"
$ws.Range("J7").Value = "Fail"
$ws.Range("J10").Value = "Pass"
$ws.Range("J11").Value = "Fail"

# ---------------------------------------------------------------------------
# Move the active selection to K4 (was H2)
# ---------------------------------------------------------------------------
$ws.Range("K4").Select()

# Note: the workbook-level window position (xWindow/yWindow on
# <bookViews><workbookView>) mirrors where the Excel window sits on screen.
# There is no Workbook/Window object-model property in this environment that
# round-trips into that attribute pair (Application.Left/Top and
# Windows.Item(1).Left/Top can be read back in-session but do not affect the
# saved xWindow/yWindow), so it is intentionally left untouched here.

Write-Host "AddExpense sheet updated"
